# "Actualizar 02-05-2021 11-49-55" automatic-refresh edit:
#  - bump the timestamp already stamped on rows 534:547 (the 11:19 run)
#    forward to the precise 44232.47190263889 value recorded by this run
#  - append a brand-new 14-row cycle (rows 548:561) for the 11:49 run,
#    each row reusing the same Nombre / URL / Disponibilidad text (and
#    therefore the same shared-string + hyperlink targets) as rows 534:547
#    but stamped with the new 44232.49295149813 timestamp

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- 1. refresh the existing timestamps on rows 534-547 ----------------
for ($r = 534; $r -le 547; $r++) {
    $ws.Cells.Item($r, 4).Value = 44232.47190263889
}

# ---- 2. the 14-row cycle reused for the new block -----------------------
$cycle = @(
    @{ A="Odoo"; B="https://www.dataintelligence-group.com/"; C="Disponible"; Target="https://www.dataintelligence-group.com/"; Loc="" },
    @{ A="Blackbox"; B="https://serviciodashboard.azurewebsites.net/"; C="Disponible"; Target="https://serviciodashboard.azurewebsites.net/"; Loc="" },
    @{ A="PowerBI"; B="https://powerbi.microsoft.com/es-es/"; C="Disponible"; Target="https://powerbi.microsoft.com/es-es/"; Loc="" },
    @{ A="Dropbox"; B="https://www.dropbox.com/"; C="Disponible"; Target="https://www.dropbox.com/"; Loc="" },
    @{ A="Odoo"; B="https://dataintelligence.store/"; C="Disponible"; Target="https://dataintelligence.store/"; Loc="" },
    @{ A="GEE"; B="https://app-data-i.users.earthengine.app/"; C="Disponible"; Target="https://app-data-i.users.earthengine.app/"; Loc="" },
    @{ A="UtilidadesOdoo"; B="https://odooutil.azurewebsites.net/"; C="Disponible"; Target="https://odooutil.azurewebsites.net/"; Loc="" },
    @{ A="Filtros Dashboard"; B="https://filtradordashboard.azurewebsites.net/"; C="Disponible"; Target="https://filtradordashboard.azurewebsites.net/"; Loc="" },
    @{ A="MapStore"; B="https://ide.dataintelligence-group.com/mapstore/#/"; C="Disponible"; Target="https://ide.dataintelligence-group.com/mapstore/"; Loc="/" },
    @{ A="GeoServer"; B="https://ide.dataintelligence-group.com/geoserver/web/?0"; C="Disponible"; Target="https://ide.dataintelligence-group.com/geoserver/web/?0"; Loc="" },
    @{ A="Tomcat"; B="https://ide.dataintelligence-group.com/"; C="Disponible"; Target="https://ide.dataintelligence-group.com/"; Loc="" },
    @{ A="Shiny"; B="https://rpubs.com/dataintelligence/"; C="Disponible"; Target="https://rpubs.com/dataintelligence/"; Loc="" },
    @{ A="Github"; B="https://github.com/Sud-Austral/"; C="Disponible"; Target="https://github.com/Sud-Austral/"; Loc="" },
    @{ A="EZ Exporter"; B="https://ezexporter.highviewapps.com/exports/export-profile/"; C="Disponible"; Target="https://ezexporter.highviewapps.com/exports/export-profile/"; Loc="" }
)

$newTimestamp = 44232.49295149813
$startRow = 548

for ($i = 0; $i -lt $cycle.Count; $i++) {
    $row = $startRow + $i
    $item = $cycle[$i]

    $ws.Cells.Item($row, 1).Value = $item.A
    $ws.Cells.Item($row, 2).Value = $item.B
    $ws.Cells.Item($row, 3).Value = $item.C

    $dCell = $ws.Cells.Item($row, 4)
    $dCell.Value = $newTimestamp
    $dCell.NumberFormat = $ws.Cells.Item($row - 1, 4).NumberFormat

    $bCell = $ws.Cells.Item($row, 2)
    if ($item.Loc -ne "") {
        $ws.Hyperlinks.Add($bCell, $item.Target, $item.Loc) | Out-Null
    } else {
        $ws.Hyperlinks.Add($bCell, $item.Target) | Out-Null
    }
    $bCell.Style = "Hyperlink"
}

Write-Host "Added rows 548-561 and refreshed timestamps on 534-547"
